$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 16-23 hold one "Periodo Mora" record each (columns E=periodo, F=valor mora, G=salario basico).
# The update re-sorts the periods into ascending order (2009 -> 2106, they previously ran
# descending 2106 -> 2009) while keeping each periodo's "valor mora" (F) value attached to it,
# and refreshes every "salario basico" (G) figure from 908526 to the new 877803.

$periodos = @("2009", "2010", "2011", "2012", "2102", "2104", "2105", "2106")
$valorMora = @(35112, 35112, 35112, 35112, 36341, 36341, 35112, 24578)
$salarioBasico = 877803

for ($i = 0; $i -lt 8; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periodos[$i]
    $ws.Range("F$row").Value = $valorMora[$i]
    $ws.Range("G$row").Value = $salarioBasico
}
